# Updated table 1 with changed missing
#
# Re-writes a set of summary-table cells on both the "unadjusted" and
# "adjusted" worksheets. Every target cell holds plain text (counts like
# "4", percentages like "7%", or formatted numbers like "1,411"), so we
# force text storage (NumberFormat "@") before writing the literal value,
# then restore the cell's original "Normal" style so no stray formatting
# is left behind.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# ---- "unadjusted" sheet ----
$ws1 = $wb.Worksheets.Item("unadjusted")

Set-TextValue $ws1 "B19" "4"
Set-TextValue $ws1 "C19" "3"
Set-TextValue $ws1 "D19" "1"

Set-TextValue $ws1 "B23" "8"
Set-TextValue $ws1 "C23" "3"
Set-TextValue $ws1 "D23" "5"

Set-TextValue $ws1 "B28" "3"
Set-TextValue $ws1 "C28" "1"
Set-TextValue $ws1 "D28" "2"

Set-TextValue $ws1 "B30" "415 (56%)"
Set-TextValue $ws1 "C30" "112 (46%)"
Set-TextValue $ws1 "D30" "303 (60%)"

Set-TextValue $ws1 "B31" "203 (27%)"
Set-TextValue $ws1 "C31" "82 (34%)"
Set-TextValue $ws1 "D31" "121 (24%)"

Set-TextValue $ws1 "B32" "104 (14%)"
Set-TextValue $ws1 "C32" "38 (16%)"
Set-TextValue $ws1 "D32" "66 (13%)"

Set-TextValue $ws1 "B33" "21 (3%)"
Set-TextValue $ws1 "C33" "9 (4%)"
Set-TextValue $ws1 "D33" "12 (2%)"

Set-TextValue $ws1 "B34" "1,411"
Set-TextValue $ws1 "C34" "267"
Set-TextValue $ws1 "D34" "1,144"

Set-TextValue $ws1 "B39" "2"
Set-TextValue $ws1 "C39" "1"
Set-TextValue $ws1 "D39" "1"

# ---- "adjusted" sheet ----
$ws2 = $wb.Worksheets.Item("adjusted")

Set-TextValue $ws2 "D16" "7%"

Set-TextValue $ws2 "B19" "71,192"
Set-TextValue $ws2 "C19" "33,682"
Set-TextValue $ws2 "D19" "37,510"

Set-TextValue $ws2 "B22" "100%"
Set-TextValue $ws2 "C22" "98%"

Set-TextValue $ws2 "B23" "292,666"
Set-TextValue $ws2 "C23" "196,914"
Set-TextValue $ws2 "D23" "95,752"

Set-TextValue $ws2 "B28" "55,579"
Set-TextValue $ws2 "C28" "9,230"
Set-TextValue $ws2 "D28" "46,348"

Set-TextValue $ws2 "B30" "62%"
Set-TextValue $ws2 "C30" "52%"
Set-TextValue $ws2 "D30" "66%"

Set-TextValue $ws2 "B31" "22%"
Set-TextValue $ws2 "C31" "28%"
Set-TextValue $ws2 "D31" "20%"

Set-TextValue $ws2 "B32" "13%"
Set-TextValue $ws2 "C32" "16%"
Set-TextValue $ws2 "D32" "11%"

Set-TextValue $ws2 "B33" "3%"
Set-TextValue $ws2 "C33" "4%"
Set-TextValue $ws2 "D33" "2%"

Set-TextValue $ws2 "B34" "48,238,432"
Set-TextValue $ws2 "C34" "7,992,154"
Set-TextValue $ws2 "D34" "40,246,278"

Set-TextValue $ws2 "C37" "73%"

Set-TextValue $ws2 "B39" "105,951"
Set-TextValue $ws2 "C39" "40,321"
Set-TextValue $ws2 "D39" "65,630"
